# Update the existing data row (row 2) with a new invoice (SONOVISION
# ELECTRONICS replaces WESTSIDE), then append a brand-new invoice row
# (row 3, LAKSHMI AGENCIES) with the same layout/formatting conventions
# used by the rest of the "GST Report" sheet.
#
# Every column on this sheet holds plain text -- including the
# numeric-/date-looking ones -- so cells whose content could otherwise be
# auto-parsed by Excel into a date or a number are briefly switched to a
# text number format before the value is written, then restored to the
# Normal style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Row 2: replace WESTSIDE invoice with SONOVISION ELECTRONICS invoice ---
$ws.Range("A2").Value = "SONOVISION ELECTRONICS`nPVT LTD"
$ws.Range("B2").Value = "37ABCCS7530B1ZK"
$ws.Range("C2").Value = "NDYL 3826"
Set-TextValue $ws.Range("D2") "17/Feb/2023"
Set-TextValue $ws.Range("E2") "69000"
Set-TextValue $ws.Range("F2") "15094"
Set-TextValue $ws.Range("G2") "7547"
Set-TextValue $ws.Range("H2") "7547"
$ws.Range("I2").Value = "N/A"
$ws.Range("J2").Value = "85287219`n0"

# --- Row 3: new invoice row for LAKSHMI AGENCIES ---
$ws.Range("A3").Value = "LAKSHMI AGENCIES`nNo:18, Kannadasan Nagar Main Road,`nRamapuram"
$ws.Range("B3").Value = "33AABFL7718B1ZQ"
$ws.Range("C3").Value = "LA226412507098"
$ws.Range("D3").Value = "17/07/2025"
Set-TextValue $ws.Range("E3") "33725.00"
Set-TextValue $ws.Range("F3") "1440.00"
Set-TextValue $ws.Range("G3") "720.00"
Set-TextValue $ws.Range("H3") "720.00"
$ws.Range("I3").Value = "N/A"
$ws.Range("J3").Value = "15121910`n15121910`n15121910`n15180039`n15180039"

# Match the wrap-text / top-vertical-align style already used for J2's
# multi-line HSN codes.
$ws.Range("J3").WrapText = $true
$ws.Range("J3").VerticalAlignment = -4160

$ws.Range("A1:J3").Select
